# Update the two-digit-divided-by-one-digit division problems in the table.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1).Range
$ok = $cell.Find.Execute("66÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷8=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (1,1): 66÷9= -> 93÷8=" }
$cell = $t.Cell(1, 2).Range
$ok = $cell.Find.Execute("80÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷4=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (1,2): 80÷3= -> 71÷4=" }
$cell = $t.Cell(1, 3).Range
$ok = $cell.Find.Execute("38÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷3=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (1,3): 38÷7= -> 84÷3=" }
$cell = $t.Cell(1, 4).Range
$ok = $cell.Find.Execute("56÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷2=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (1,4): 56÷3= -> 20÷2=" }
$cell = $t.Cell(1, 5).Range
$ok = $cell.Find.Execute("55÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷6=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (1,5): 55÷5= -> 26÷6=" }
$cell = $t.Cell(5, 1).Range
$ok = $cell.Find.Execute("15÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (5,1): 15÷4= -> 18÷3=" }
$cell = $t.Cell(5, 2).Range
$ok = $cell.Find.Execute("82÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷3=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (5,2): 82÷3= -> 51÷3=" }
$cell = $t.Cell(5, 3).Range
$ok = $cell.Find.Execute("25÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷8=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (5,3): 25÷2= -> 12÷8=" }
$cell = $t.Cell(5, 4).Range
$ok = $cell.Find.Execute("78÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷9=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (5,4): 78÷2= -> 79÷9=" }
$cell = $t.Cell(5, 5).Range
$ok = $cell.Find.Execute("24÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "54÷3=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (5,5): 24÷7= -> 54÷3=" }
$cell = $t.Cell(9, 1).Range
$ok = $cell.Find.Execute("82÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷6=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (9,1): 82÷2= -> 29÷6=" }
$cell = $t.Cell(9, 2).Range
$ok = $cell.Find.Execute("59÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷3=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (9,2): 59÷6= -> 42÷3=" }
$cell = $t.Cell(9, 3).Range
$ok = $cell.Find.Execute("80÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷4=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (9,3): 80÷8= -> 53÷4=" }
$cell = $t.Cell(9, 4).Range
$ok = $cell.Find.Execute("57÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷4=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (9,4): 57÷7= -> 47÷4=" }
$cell = $t.Cell(9, 5).Range
$ok = $cell.Find.Execute("90÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷8=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (9,5): 90÷9= -> 88÷8=" }
$cell = $t.Cell(13, 1).Range
$ok = $cell.Find.Execute("97÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷8=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (13,1): 97÷7= -> 62÷8=" }
$cell = $t.Cell(13, 2).Range
$ok = $cell.Find.Execute("42÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷5=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (13,2): 42÷5= -> 55÷5=" }
$cell = $t.Cell(13, 3).Range
$ok = $cell.Find.Execute("23÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (13,3): 23÷2= -> 93÷7=" }
$cell = $t.Cell(13, 4).Range
$ok = $cell.Find.Execute("61÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷3=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (13,4): 61÷8= -> 55÷3=" }
$cell = $t.Cell(13, 5).Range
$ok = $cell.Find.Execute("20÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "96÷2=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (13,5): 20÷9= -> 96÷2=" }
$cell = $t.Cell(17, 1).Range
$ok = $cell.Find.Execute("34÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷8=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (17,1): 34÷2= -> 37÷8=" }
$cell = $t.Cell(17, 2).Range
$ok = $cell.Find.Execute("18÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷7=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (17,2): 18÷9= -> 98÷7=" }
$cell = $t.Cell(17, 3).Range
$ok = $cell.Find.Execute("71÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷3=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (17,3): 71÷7= -> 14÷3=" }
$cell = $t.Cell(17, 4).Range
$ok = $cell.Find.Execute("71÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷6=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (17,4): 71÷5= -> 42÷6=" }
$cell = $t.Cell(17, 5).Range
$ok = $cell.Find.Execute("22÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷6=", 2)
if (-not $ok) { Write-Host "WARNING: replace failed for cell (17,5): 22÷4= -> 87÷6=" }
